# Map each changed cell (A1 ref) to its new text value.
# These columns (Price / Volume(1h)%) are stored as plain TEXT in the
# sheet, not numbers, so we briefly force a text NumberFormat while
# writing the value (otherwise Excel auto-coerces a numeric-looking
# string like "328.47" or "-0.66%" into a real number/percentage) and
# then restore NumberFormat back to General so the cell format is left
# exactly as it was found.
$updates = [ordered]@{
    "D2" = "328.47"
    "E2" = "-0.66%"
    "D3" = "43.76"
    "E3" = "5.59%"
    "D4" = "5.436"
    "E4" = "-3.87%"
    "D5" = "0.08148"
    "E5" = "-2.30%"
    "D6" = "8.730"
    "E6" = "-0.59%"
    "D7" = "1.912"
    "E7" = "-3.24%"
    "D8" = "4.318"
    "E8" = "-3.61%"
    "D10" = "0.9425"
    "E10" = "1.93%"
    "D11" = "0.1178"
    "E11" = "-8.84%"
    "D12" = "0.1895"
    "E12" = "-3.16%"
    "D13" = "0.09759"
    "E13" = "3.67%"
    "D14" = "0.04249"
    "E14" = "8.76%"
    "D15" = "0.1068"
    "E15" = "1.13%"
    "D16" = "0.001279"
    "E16" = "-1.54%"
    "D17" = "0.006042"
    "E17" = "-0.96%"
    "D18" = "3.551"
    "E18" = "3.23%"
    "E19" = "-0.53%"
    "D20" = "8.789"
    "E20" = "-2.84%"
    "D21" = "0.1358"
    "E21" = "-0.87%"
    "D23" = "0.04379"
    "E23" = "-0.67%"
    "D24" = "0.001242"
    "E24" = "-2.46%"
    "D25" = "0.004351"
    "E25" = "-1.03%"
    "D26" = "0.0001237"
    "E26" = "3.04%"
    "D27" = "0.0004007"
    "E27" = "31.58%"
    "D39" = "0.02675"
    "E39" = "-5.65%"
    "D40" = "0.05630"
    "E40" = "2.14%"
    "D41" = "0.007888"
    "E41" = "1.35%"
    "D42" = "0.009779"
    "E42" = "5.09%"
    "D43" = "0.1407"
    "E43" = "-1.96%"
    "D44" = "0.002129"
    "E44" = "-0.51%"
    "D45" = "0.009624"
    "E45" = "-12.73%"
    "D46" = "0.00007037"
    "E46" = "-0.70%"
    "E47" = "0.45%"
    "D48" = "0.003468"
    "E48" = "-2.23%"
    "E49" = "0.18%"
    "E50" = "0.45%"
    "E51" = "0.45%"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.NumberFormat = "General"
}
